# Update currency year to 2023
#
# The workbook's "About" sheet holds a currency-conversion factor in A26
# together with a descriptive label in B26 ("<year> dollars per 2012 dollar").
# This commit refreshes the conversion factor and its label from the 2021
# vintage to the 2023 vintage. Dependent sheets (OCCF-DpLOCU, OCCF-DpMOCU,
# OCCF-DpSOCU) hold formulas that reference About!A26, so they recompute
# automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

# New 2012->2023 dollar conversion factor and its label.
$ws.Range("A26").Value = 0.75350342301658668
$ws.Range("B26").Value = "2023 dollars per 2012 dollar"

# Reflect the author's last on-sheet selection (cell A26) when the file was saved.
$ws.Activate() | Out-Null
$ws.Range("A26").Select() | Out-Null
